# Commit: "update the scrape file and update links in data, progresson nlp_label"
#
# Net effect on the worksheet:
#   - Cell B11 ("Innovation on Board" row) gets refreshed link text: the old
#     2017 FAQ link is replaced with new 2024 schedule/FAQ links.
#   - Row 11 grows taller (51 -> 68) to fit the extra wrapped line.
#   - The view had scrolled/selected down near the bottom (A12/B21); it's
#     reset back up near the edited row (A8/B11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "Innovation on Board" link cell (B11) ---
$newLink = "https://www.innovationboard.ca/ | " + [char]10 + `
    "https://www.innovationboard.ca//learn-more.html | " + [char]10 + `
    "https://www.innovationboard.ca/2024/06/23/schedule.html | https://www.innovationboard.ca/2024/06/23/faq.html"

$ws.Range("B11").Value = $newLink

# Row grew to fit the new (longer) wrapped content
$ws.Rows.Item(11).RowHeight = 68

# --- Update the window/sheet view back to the edited area ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()

# Window had scrolled horizontally on screen as well
$excel.ActiveWindow.Left = 15880
